$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.180173397064209
$ws.Range("B1").Value = 2.164919376373291
$ws.Range("C1").Value = 3.625627040863037
$ws.Range("D1").Value = 3.418977975845337
$ws.Range("E1").Value = 1.152811288833618
